$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.805.94"
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").Value = "1.888.92"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "0.7692"
$ws.Range("E5").Value = "  -5.25%  "
$ws.Range("D6").Value = "244.28"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.3120"
$ws.Range("E8").Value = "  -4.34%  "
$ws.Range("E9").Value = "  -7.12%  "
$ws.Range("D10").Value = "0.07218"
$ws.Range("E10").Value = "  -0.38%  "
$ws.Range("D11").Value = "0.08096"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").Value = "0.7638"
$ws.Range("E12").Value = "  -3.66%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.938.46"
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.494"
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("D15").Value = "92.25"
$ws.Range("E15").Value = "  -2.05%  "
$ws.Range("D16").Value = "6.141"
$ws.Range("E16").Value = "  +0.94%  "
$ws.Range("D17").Value = "29.817.33"
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("D18").Value = "13.89"
$ws.Range("E18").Value = "  -2.88%  "
$ws.Range("D19").Value = "242.70"
$ws.Range("E19").Value = "  -3.33%  "
$ws.Range("D20").Value = "0.000007765"
$ws.Range("E20").Value = "  -1.16%  "
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.160.27"
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "8.190"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "0.1564"
$ws.Range("E25").Value = "  -6.51%  "
$ws.Range("D26").Value = "9.433"
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").Value = "162.29"
$ws.Range("E27").Value = "  -2.98%  "
$ws.Range("D28").Value = "18.75"
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("D29").Value = "2.036"
$ws.Range("E29").Value = "  -5.71%  "
$ws.Range("D30").Value = "1.463"
$ws.Range("E30").Value = "  +6.53%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").Value = "4.447"
$ws.Range("E32").Value = "  +2.22%  "
$ws.Range("D33").Value = "4.075"
$ws.Range("E33").Value = "  -1.55%  "
$ws.Range("D34").Value = "0.05520"
$ws.Range("E34").Value = "  -2.86%  "
$ws.Range("D35").Value = "1.255"
$ws.Range("E35").Value = "  -3.55%  "
$ws.Range("D36").Value = "0.7474"
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("D38").Value = "2.635"
$ws.Range("E38").Value = "  -3.44%  "
$ws.Range("D39").Value = "0.01919"
$ws.Range("E39").Value = "  -2.06%  "
$ws.Range("E40").Value = "  -1.37%  "
$ws.Range("D41").Value = "1.142.57"
$ws.Range("E41").Value = "  +9.51%  "
$ws.Range("D42").Value = "73.67"
$ws.Range("E42").Value = "  -1.32%  "
$ws.Range("E43").Value = "  -1.91%  "
$ws.Range("D44").Value = "5.904"
$ws.Range("E44").Value = "  -1.42%  "
$ws.Range("D45").Value = "0.8470"
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "102.81"
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("D48").Value = "1.881"
$ws.Range("E48").Value = "  -2.44%  "
$ws.Range("D49").Value = "9.892"
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "7.431"
$ws.Range("E50").Value = "  -2.71%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").Value = "3.012"
$ws.Range("E51").Value = "  -2.32%  "
